$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original style of column D (Price) data range, then force text
# format while assigning numeric-looking strings so Excel keeps them as text
# (matching the source inlineStr cells), then restore the original style so
# no stray style/number-format artifacts are introduced.
$dRange = $ws.Range("D2:D51")
$dStyle = $dRange.Style
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.664.97"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "1.631.42"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "213.46"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +3.67%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("E11").Value = "  +3.30%  "
$ws.Range("D12").Value = "1.857.89"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.618.30"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "26.657.86"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "63.50"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("D19").Value = "218.70"
$ws.Range("E19").Value = "  +8.04%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "4.30"
$ws.Range("E21").Value = "  +0.85%  "
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +1.83%  "
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("E24").Value = "  +4.60%  "
$ws.Range("D25").Value = "147.72"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("D29").Value = "15.52"
$ws.Range("E29").Value = "  +2.11%  "
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  -2.57%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E32").Value = "  +3.95%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  +0.90%  "
$ws.Range("D36").Value = "1.214.01"
$ws.Range("E36").Value = "  +4.84%  "
$ws.Range("E37").Value = "  +4.56%  "
$ws.Range("D38").Value = "0.806"
$ws.Range("E38").Value = "  +0.25%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("D40").Value = "0.501"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  -1.95%  "
$ws.Range("D42").Value = "0.794"
$ws.Range("E42").Value = "  +1.33%  "
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "1.768.89"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "92.68"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("D46").Value = "1.55"
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "55.07"
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0104"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").Value = "7.62"
$ws.Range("E50").Value = "  +4.10%  "
$ws.Range("D51").Value = "0.409"
$ws.Range("E51").Value = "  -0.10%  "

# Restore original style/number format for column D
$dRange.Style = $dStyle
